# Updates documentacion/peticiones.xlsx to reflect progress through call 14
# ("Terminada hasta la llamada 14, postman actualizado, excel actualizado").
#
# Changes:
#  - Rows 9-14 ("terminado" column I): "no" -> "si"
#  - F11 body sample: "id: id," -> "id: posicion_de_la_instancia,"
#  - F13 body sample: expanded to the new atributo/valores shape
#  - H13 comments: add a note to check Postman for an example of this call
#  - Selection / top-left cell moved to reflect the new scroll position (A13 / I15)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "terminado" column now marked "si" for the finished calls (rows 9-14)
$ws.Range("I9").Value = "si"
$ws.Range("I10").Value = "si"
$ws.Range("I11").Value = "si"
$ws.Range("I12").Value = "si"
$ws.Range("I13").Value = "si"
$ws.Range("I14").Value = "si"

# "modificar instancia" body sample now references the row position instead of a generic id
$ws.Range("F11").Value = "{`nid: posicion_de_la_instancia,`nobjeto: objeto`n}"

# "agregar atributo" body sample expanded with the full attribute configuration/values shape
$ws.Range("F13").Value = "{`natributo: { `nobjeto_configuracion_atributo`n},`nvalores:  [{ `nnombre_atributo: valor `n}]"

# New comment directing to Postman for an example of this call
$ws.Range("H13").Value = "Revisar postman para ver el ejemplo de esta llamada"

# Scroll/selection moved on to the next pending call
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("I15").Select()
